$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 180.5625
$ws.Range("I2").Value = 176.66667
$ws.Range("J2").Value = 185.57143
$ws.Range("K2").Value = 176.66667
$ws.Range("L2").Value = 185.57143
$ws.Range("M2").Value = -63.66667000000001
$ws.Range("N2").Value = -411.57143
$ws.Range("H127").Value = 928.17645
$ws.Range("I127").Value = 338.66666
$ws.Range("J127").Value = 1140.4
$ws.Range("K127").Value = 1015.99998
$ws.Range("L127").Value = 3421.2
$ws.Range("M127").Value = 3944.00002
$ws.Range("N127").Value = -13341.2
$ws.Range("H129").Value = 988.84
$ws.Range("I129").Value = 434.81818
$ws.Range("J129").Value = 1057.3146
$ws.Range("K129").Value = 1304.45454
$ws.Range("L129").Value = 3171.9438
$ws.Range("M129").Value = 3695.54546
$ws.Range("N129").Value = -13171.9438

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 936.6316
$ws.Range("I110").Value = 690.5333000000001
$ws.Range("K110").Value = 690.5333000000001
$ws.Range("M110").Value = 1354.4667
$ws.Range("H122").Value = 1312.3889
$ws.Range("I122").Value = 1338.875
$ws.Range("J122").Value = 1100.5
$ws.Range("K122").Value = 4016.625
$ws.Range("L122").Value = 3301.5
$ws.Range("M122").Value = -1566.625
$ws.Range("N122").Value = -8201.5
$ws.Range("H138").Value = 20429
$ws.Range("J138").Value = 20429
$ws.Range("L138").Value = 20429
$ws.Range("N138").Value = -30709

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1744.5
$ws.Range("I99").Value = 1049.1666
$ws.Range("J99").Value = 3830.5
$ws.Range("K99").Value = 1049.1666
$ws.Range("L99").Value = 3830.5
$ws.Range("M99").Value = 448.8334
$ws.Range("N99").Value = -6826.5
$ws.Range("H105").Value = 2766.847
$ws.Range("I105").Value = 1776.125
$ws.Range("J105").Value = 2854.9111
$ws.Range("K105").Value = 1776.125
$ws.Range("L105").Value = 2854.9111
$ws.Range("M105").Value = -29.125
$ws.Range("N105").Value = -6348.911099999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3400.3333
$ws.Range("I16").Value = 1250.2
$ws.Range("J16").Value = 4936.143
$ws.Range("K16").Value = 1250.2
$ws.Range("L16").Value = 4936.143
$ws.Range("M16").Value = -963.2
$ws.Range("N16").Value = -5510.143
$ws.Range("H105").Value = 1155
$ws.Range("I105").Value = 1225
$ws.Range("J105").Value = 1050
$ws.Range("K105").Value = 1225
$ws.Range("L105").Value = 1050
$ws.Range("M105").Value = 522
$ws.Range("N105").Value = -4544
$ws.Range("H113").Value = 3400.3333
$ws.Range("I113").Value = 1250.2
$ws.Range("J113").Value = 4936.143
$ws.Range("K113").Value = 1250.2
$ws.Range("L113").Value = 4936.143
$ws.Range("M113").Value = 919.8
$ws.Range("N113").Value = -9276.143
$ws.Range("H132").Value = 2833.5833
$ws.Range("I132").Value = 1991.9333
$ws.Range("K132").Value = 5975.7999
$ws.Range("M132").Value = -3445.7999
$ws.Range("H134").Value = 2333.4285
$ws.Range("I134").Value = 1453.9546
$ws.Range("K134").Value = 4361.8638
$ws.Range("M134").Value = -1826.8638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 442.44446
$ws.Range("I92").Value = 470.66666
$ws.Range("J92").Value = 428.33334
$ws.Range("K92").Value = 1411.99998
$ws.Range("L92").Value = 1285.00002
$ws.Range("M92").Value = -163.9999800000001
$ws.Range("N92").Value = -3781.00002
$ws.Range("H93").Value = 4732.6665
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 4732.6665
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 14197.9995
$ws.Range("M93").Value = $null
$ws.Range("N93").Value = -17941.9995
$ws.Range("H95").Value = 4000
$ws.Range("I95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("M95").Value = $null
$ws.Range("H96").Value = 34413.637
$ws.Range("J96").Value = 34413.637
$ws.Range("L96").Value = 103240.911
$ws.Range("N96").Value = -107358.911
$ws.Range("H113").Value = 565815.3
$ws.Range("I113").Value = 1231987.6
$ws.Range("J113").Value = 578.2121
$ws.Range("K113").Value = 3695962.8
$ws.Range("L113").Value = 1734.6363
$ws.Range("M113").Value = -3693792.8
$ws.Range("N113").Value = -6074.6363

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6109.3447
$ws.Range("I70").Value = 6191.5186
$ws.Range("K70").Value = 6191.5186
$ws.Range("M70").Value = -5921.5186
$ws.Range("H73").Value = 6109.3447
$ws.Range("I73").Value = 6191.5186
$ws.Range("K73").Value = 6191.5186
$ws.Range("M73").Value = -5255.5186
$ws.Range("H122").Value = 3449810.8
$ws.Range("I122").Value = 4763263
$ws.Range("K122").Value = 14289789
$ws.Range("M122").Value = -14287339
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2133.3333
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 2142.8572
$ws.Range("K68").Value = 2000
$ws.Range("L68").Value = 2142.8572
$ws.Range("M68").Value = -1251
$ws.Range("N68").Value = -3640.8572
$ws.Range("H71").Value = 2133.3333
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 2142.8572
$ws.Range("K71").Value = 10000
$ws.Range("L71").Value = 10714.286
$ws.Range("M71").Value = -6256
$ws.Range("N71").Value = -18202.286
$ws.Range("H111").Value = 47900
$ws.Range("J111").Value = 47900
$ws.Range("L111").Value = 47900
$ws.Range("N111").Value = -56080
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").Value = $null
$ws.Range("H122").Value = 2163.524
$ws.Range("I122").Value = 2083.375
$ws.Range("J122").Value = 2420
$ws.Range("K122").Value = 6250.125
$ws.Range("L122").Value = 7260
$ws.Range("M122").Value = -3800.125
$ws.Range("N122").Value = -12160
$ws.Range("H136").Value = 1641.7561
$ws.Range("I136").Value = 1272.4839
$ws.Range("J136").Value = 2786.5
$ws.Range("K136").Value = 3817.4517
$ws.Range("L136").Value = 8359.5
$ws.Range("M136").Value = -1267.4517
$ws.Range("N136").Value = -13459.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2270.2563
$ws.Range("I122").Value = 2083.6072
$ws.Range("J122").Value = 2745.3635
$ws.Range("K122").Value = 6250.821599999999
$ws.Range("L122").Value = 8236.0905
$ws.Range("M122").Value = -3800.821599999999
$ws.Range("N122").Value = -13136.0905
$ws.Range("H141").Value = 230000
$ws.Range("J141").Value = 230000
$ws.Range("L141").Value = 230000
$ws.Range("N141").Value = -240360
